# Swap the match-record data (columns B:AB) between each of the following
# row pairs on the "Montenegro Prva Liga" sheet. Column A (the running
# index number) stays put on its own row; every other field (match id,
# teams, odds, etc.) for the two rows trades places, which is what the
# target diff shows for all eight pairs below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(48, 49),
    @(75, 76),
    @(100, 101),
    @(121, 122),
    @(126, 127),
    @(152, 153),
    @(154, 155),
    @(160, 161)
)

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

foreach ($pair in $rowPairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    foreach ($col in $cols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"

        $val1 = $ws.Range($addr1).Value2
        $val2 = $ws.Range($addr2).Value2

        # Only touch cells whose value actually changes, to avoid needless
        # re-writes of untouched data.
        if ($val1 -ne $val2) {
            $ws.Range($addr1).Value2 = $val2
            $ws.Range($addr2).Value2 = $val1
        }
    }
}
